# Weekly refresh: insert a new (most-recent) week of "Perejil" price rows at
# the top of the data block (rows 126/127), pushing the existing history
# down by one week (2 rows: Primera/Segunda), and append the resulting
# surplus pair at the bottom (rows 154/155).
#
# Net effect: the data that used to live in rows 126..153 now lives in
# rows 128..155, a brand-new pair of rows (126/127) is created with an
# updated date, and the dimension grows from R153 to R155.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing weekly block (rows 126-153) down by inserting two new
# blank rows at 126-127. Excel's Insert pushes rows 126-153 -> 128-155 and
# grows the used range / dimension automatically.
$ws.Rows("126:127").Insert()

# The row that used to be 126 (Primera) is now at 128, and the old 127
# (Segunda) is now at 129. Duplicate that pair's full contents into the new
# 126/127 rows, then overwrite just the date (column D) with the new week.
$ws.Range("A126:R126").Value = $ws.Range("A128:R128").Value()
$ws.Range("A127:R127").Value = $ws.Range("A129:R129").Value()

$ws.Cells.Item(126, 4).Value = 44754
$ws.Cells.Item(127, 4).Value = 44754
